$d = $word.ActiveDocument

$wNs   = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$pkgNs = "http://schemas.microsoft.com/office/2006/xmlPackage"

function New-PackageXml([string]$bodyXml) {
    return "<pkg:package xmlns:pkg=`"$pkgNs`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"$wNs`"><w:body>$bodyXml</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
}

# ---------------------------------------------------------------------------
# 1) Elaborate on the ETL-methods sentence: explain the languages used and the
#    public data sources the pipeline draws on.
# ---------------------------------------------------------------------------
$oldParaText = "The methods described in the paper use free, open source software as well as " + `
    "Extract Transform and Load (ETL) methods written by the author to computationally discover " + `
    "new information using publicly available data.  It is this computational  content which " + `
    "makes it suitable for inclusion in this journal."

$rngEtl = $d.Content
$found = $rngEtl.Find.Execute($oldParaText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find ETL paragraph text" }
$rngEtl.Expand(4) | Out-Null

$etlBody = '<w:p>' + `
    '<w:r><w:t xml:space="preserve">The methods described in the paper use free, open source software as well as </w:t></w:r>' + `
    '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Extract Transform and Load </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">(ETL) methods written by the author </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">using Perl, bash, and R </w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:t>to computationally discover new information using publicly available data</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> from </w:t></w:r>' + `
    '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">myChEMBL, NCBI, and </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:i/></w:rPr><w:t>EnsemblGenomes</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>.</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">  It is this </w:t></w:r>' + `
    '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">computational </w:t></w:r>' + `
    '<w:r><w:t>content which makes it suitable for inclusion in this journal.</w:t></w:r>' + `
    '</w:p>'

$rngEtl.InsertXML((New-PackageXml $etlBody))

# ---------------------------------------------------------------------------
# 2) Bio paragraph: merge the "Senior I.T. developer..." runs into one plain
#    run (dropping the spell-check markers around Clarivate) and drop the
#    stray _GoBack bookmark that used to sit here (it has moved above).
# ---------------------------------------------------------------------------
$oldBioText = "Senior I.T. developer at Clarivate Analytics (Formerly Thomson Reuters Life Sciences.)"

$rngBio = $d.Content
$found = $rngBio.Find.Execute($oldBioText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find bio paragraph text" }
$rngBio.Expand(4) | Out-Null

$bioBody = '<w:p><w:r><w:t>Senior I.T. developer at Clarivate Analytics (Formerly Thomson Reuters Life Sciences.)</w:t></w:r></w:p>'
$rngBio.InsertXML((New-PackageXml $bioBody))

# InsertXML re-creates a _GoBack bookmark around whatever range it just
# touched (mirroring Word's own "remember last edit" behaviour) because this
# paragraph used to own that bookmark. It belongs in the ETL paragraph now,
# so remove this stray copy.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 3) Signature block: mark the page break that now falls right before the
#    "e-mail:" line once the letter grew by a paragraph.
# ---------------------------------------------------------------------------
$oldMailText = "e-mail:" + [char]9 + "jsinger@rcn.com "

$rngMail = $d.Content
$found = $rngMail.Find.Execute($oldMailText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find e-mail paragraph text" }
$rngMail.Expand(4) | Out-Null

$mailBody = '<w:p>' + `
    '<w:r><w:rPr><w:iCs/></w:rPr><w:lastRenderedPageBreak/><w:t>e-mail:</w:t></w:r>' + `
    '<w:r><w:rPr><w:iCs/></w:rPr><w:tab/><w:t xml:space="preserve">jsinger@rcn.com </w:t></w:r>' + `
    '</w:p>'
$rngMail.InsertXML((New-PackageXml $mailBody))

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

Write-Host "Edits applied"
